$wb = $excel.ActiveWorkbook

# Sheet1: A1 value change ("26561,21" -> "6e-05").
# Force the cell format to Text first so the scientific-notation-looking
# string is stored verbatim instead of being parsed into a number.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").NumberFormat = "@"
$ws1.Range("A1").Value = "6e-05"

# data sheet: update suggestion coin name + sheet/cell reference
$ws2 = $wb.Worksheets.Item("data")
$ws2.Range("B1").Value = "terra-luna"
$ws2.Range("B2").Value = "Sheet1"
$ws2.Range("B3").Value = "A1"
